$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3004949.14
$ws.Range("C9").Value = 477251.99
$ws.Range("D9").Value = 3482201.13
$ws.Range("E9").Value = 13.70546881650113
$ws.Range("F9").Value = 86.29453118349888
$ws.Range("G9").Value = -53.87597485462972
$ws.Range("H9").Value = -45.73477209134096
$ws.Range("I9").Value = 29911
$ws.Range("J9").Value = 1283
$ws.Range("K9").Value = 31194
$ws.Range("L9").Value = 21531
$ws.Range("M9").Value = 161.7296516650411
$ws.Range("N9").Value = 10.41612725750623
